# Generate Report for Handoff
# Adds a new row (row 3) to each of the three sheets (Overview, zh-cn, de-de)
# representing a newly-handed-off file "9eae8c99-...".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$fileName    = "9eae8c99-54d0-400d-9277-be7a534b01baooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$fileNameE2e = "e2e\9eae8c99-54d0-400d-9277-be7a534b01baooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$zhXlf       = "9eae8c99-54d0-400d-9277-be7a534b01baoooooooooooooooooooooooooooooooooooooooo.4af193e78dc017866f7ecc219a27d13c7ef90ae4.zh-cn.xlf"
$deXlf       = "9eae8c99-54d0-400d-9277-be7a534b01baoooooooooooooooooooooooooooooooooooooooo.4af193e78dc017866f7ecc219a27d13c7ef90ae4.de-de.xlf"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31366a5dcb09ccee727b9c9d876e3cca73b9e237/e2e/9eae8c99-54d0-400d-9277-be7a534b01baooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$ws1.Range("A3").Value = $fileName
$ws1.Range("B3").Value = $fileNameE2e
$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-27 12:25:57"

$ws1.Range("A3").Style = $ws1.Range("A2").Style
$ws1.Range("C3").Style = $ws1.Range("C2").Style
$ws1.Range("D3").Style = $ws1.Range("D2").Style
$ws1.Range("E3").Style = $ws1.Range("E2").Style
$ws1.Range("F3").Style = $ws1.Range("F2").Style
$ws1.Range("G3").Style = $ws1.Range("G2").Style

$ws1.Hyperlinks.Add($ws1.Range("B3"), $hyperlinkUrl, "", "", $fileNameE2e)
$ws1.Range("B3").Font.Underline = $ws1.Range("B2").Font.Underline
$ws1.Range("B3").Font.Color = $ws1.Range("B2").Font.Color

$ws1.Columns.Item(5).ColumnWidth = 17.2159881591797
$ws1.Columns.Item(6).ColumnWidth = 17.2159881591797

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$ws2.Range("A3").Value = $fileName
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "False"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = "2016-08-27 12:25:53"
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("M3").Value = "True"
$ws2.Range("O3").Value = "False"

$ws2.Range("B3").Style = $ws2.Range("B2").Style
$ws2.Range("C3").Style = $ws2.Range("C2").Style
$ws2.Range("D3").Style = $ws2.Range("D2").Style
$ws2.Range("E3").Style = $ws2.Range("E2").Style
$ws2.Range("F3").Style = $ws2.Range("F2").Style
$ws2.Range("G3").Style = $ws2.Range("G2").Style
$ws2.Range("H3").Style = $ws2.Range("H2").Style
$ws2.Range("K3").Style = $ws2.Range("K2").Style
$ws2.Range("M3").Style = $ws2.Range("M2").Style
$ws2.Range("O3").Style = $ws2.Range("O2").Style

$ws2.Hyperlinks.Add($ws2.Range("A3"), $hyperlinkUrl, "", "", $fileName)
$ws2.Range("A3").Font.Underline = $ws2.Range("A2").Font.Underline
$ws2.Range("A3").Font.Color = $ws2.Range("A2").Font.Color

$ws2.Columns.Item(3).ColumnWidth = 17.2159881591797

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$ws3.Range("A3").Value = $fileName
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "False"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = "2016-08-27 12:25:57"
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("M3").Value = "True"
$ws3.Range("O3").Value = "False"

$ws3.Range("B3").Style = $ws3.Range("B2").Style
$ws3.Range("C3").Style = $ws3.Range("C2").Style
$ws3.Range("D3").Style = $ws3.Range("D2").Style
$ws3.Range("E3").Style = $ws3.Range("E2").Style
$ws3.Range("F3").Style = $ws3.Range("F2").Style
$ws3.Range("G3").Style = $ws3.Range("G2").Style
$ws3.Range("H3").Style = $ws3.Range("H2").Style
$ws3.Range("K3").Style = $ws3.Range("K2").Style
$ws3.Range("M3").Style = $ws3.Range("M2").Style
$ws3.Range("O3").Style = $ws3.Range("O2").Style

$ws3.Hyperlinks.Add($ws3.Range("A3"), $hyperlinkUrl, "", "", $fileName)
$ws3.Range("A3").Font.Underline = $ws3.Range("A2").Font.Underline
$ws3.Range("A3").Font.Color = $ws3.Range("A2").Font.Color

$ws3.Columns.Item(3).ColumnWidth = 17.2159881591797

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P3"))
